# Working on sequence recorder
# Append two new rows of translation text entries to the "Translation" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B10").Value = "SingleUseId8"
$ws.Range("C10").Value = "Default"
$ws.Range("D10").Value = "Center"
$ws.Range("E10").Value = "LTR"
$ws.Range("F10").Value = "<value>"

$ws.Range("B11").Value = "SingleUseId9"
$ws.Range("C11").Value = "Default"
$ws.Range("D11").Value = "Left"
$ws.Range("E11").Value = "LTR"
$ws.Range("F11").Value = "TEST"
